$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 109-111, shifting existing rows 109-141 down to 112-144
$ws.Rows("109:111").Insert()

# Row 109
$ws.Cells.Item(109, 1).Value = 9
$ws.Cells.Item(109, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(109, 3).Value = 'Metropolitana'
$ws.Cells.Item(109, 4).Value = 44524
$ws.Cells.Item(109, 5).Value = 13
$ws.Cells.Item(109, 6).Value = 100112003
$ws.Cells.Item(109, 7).Value = 'Ajo'
$ws.Cells.Item(109, 8).Value = 'Rosado'
$ws.Cells.Item(109, 9).Value = '1a nueva(o)'
$ws.Cells.Item(109, 10).Value = 3200
$ws.Cells.Item(109, 11).Value = 2000
$ws.Cells.Item(109, 12).Value = 2200
$ws.Cells.Item(109, 13).Value = 2100
$ws.Cells.Item(109, 14).Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Cells.Item(109, 15).Value = 'Provincia de Talagante'
$ws.Cells.Item(109, 16).Value = 105
$ws.Cells.Item(109, 17).Value = 20
$ws.Cells.Item(109, 18).Value = 'Hortaliza'

# Row 110
$ws.Cells.Item(110, 1).Value = 9
$ws.Cells.Item(110, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(110, 3).Value = 'Metropolitana'
$ws.Cells.Item(110, 4).Value = 44524
$ws.Cells.Item(110, 5).Value = 13
$ws.Cells.Item(110, 6).Value = 100112003
$ws.Cells.Item(110, 7).Value = 'Ajo'
$ws.Cells.Item(110, 8).Value = 'Rosado'
$ws.Cells.Item(110, 9).Value = '2a nueva(o)'
$ws.Cells.Item(110, 10).Value = 600
$ws.Cells.Item(110, 11).Value = 1600
$ws.Cells.Item(110, 12).Value = 1800
$ws.Cells.Item(110, 13).Value = 1700
$ws.Cells.Item(110, 14).Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Cells.Item(110, 15).Value = 'Provincia de Talagante'
$ws.Cells.Item(110, 16).Value = 85
$ws.Cells.Item(110, 17).Value = 20
$ws.Cells.Item(110, 18).Value = 'Hortaliza'

# Row 111
$ws.Cells.Item(111, 1).Value = 9
$ws.Cells.Item(111, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(111, 3).Value = 'Metropolitana'
$ws.Cells.Item(111, 4).Value = 44524
$ws.Cells.Item(111, 5).Value = 13
$ws.Cells.Item(111, 6).Value = 100112003
$ws.Cells.Item(111, 7).Value = 'Ajo'
$ws.Cells.Item(111, 8).Value = 'Rosado'
$ws.Cells.Item(111, 9).Value = 'Extra nueva (o)'
$ws.Cells.Item(111, 10).Value = 1400
$ws.Cells.Item(111, 11).Value = 2400
$ws.Cells.Item(111, 12).Value = 2600
$ws.Cells.Item(111, 13).Value = 2500
$ws.Cells.Item(111, 14).Value = '$/paquete 20 unidades (volumen en unidades)'
$ws.Cells.Item(111, 15).Value = 'Provincia de Talagante'
$ws.Cells.Item(111, 16).Value = 125
$ws.Cells.Item(111, 17).Value = 20
$ws.Cells.Item(111, 18).Value = 'Hortaliza'
